# "title for kable added more variables"
# - Renames activity "hype bdm" (row 11) to "beach day monday cali version"
# - Fixes "garry fiver &dinner" (row 19) to "garry fiver & dinner"
# - Moves the active selection to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "garry fiver & dinner" label first so the shared-string table
# ends up with it before the newly introduced "beach day monday cali version"
# entry (matches original authoring order).
$ws.Range("A19").Value = "garry fiver & dinner"
$ws.Range("A11").Value = "beach day monday cali version"

# Move the selection to D11, as in the saved workbook.
$ws.Range("D11").Select()
